# Applies the cryptocurrency price/volume refresh described in the commit
# "Updated cryptos list on Thu Apr 11 20:52:03 UTC 2024 with GitHub Actions".
#
# Rows 34-37 additionally had their coin order swapped (NEARProtocol/Mantle
# and dogwifhat/Hedera traded places), so those rows update B (name), C (link),
# D (price) and E (1h volume change); all other changed rows only update the
# D (price) and/or E (1h volume change) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimal numbers (e.g. 607.09) which Excel would
# otherwise auto-convert to a numeric cell. Force those specific cells to stay
# as text (matching the source workbook, where every data cell is a string)
# by pre-formatting them as Text before assigning the value.
$textPriceCells = @(
    'D5', 'D6', 'D7', 'D11', 'D12', 'D13', 'D14', 'D17', 'D21', 'D23', 'D24', 'D25', 'D29', 'D30', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D46', 'D48', 'D50'
)
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$updates = @{
    'D2' = '70.540.28'
    'E2' = '  +1.07%  '
    'D3' = '3.519.96'
    'E3' = '  +0.18%  '
    'E4' = '  -0.08%  '
    'D5' = '607.09'
    'E5' = '  +0.39%  '
    'D6' = '174.68'
    'E6' = '  +1.92%  '
    'D7' = '0.614'
    'E7' = '  -0.32%  '
    'D8' = '3.515.34'
    'E8' = '  +0.16%  '
    'E9' = '  -0.02%  '
    'E10' = '  -1.00%  '
    'D11' = '7.30'
    'E11' = '  +9.50%  '
    'D12' = '0.589'
    'E12' = '  +1.10%  '
    'D13' = '46.40'
    'E13' = '  -1.74%  '
    'D14' = '0.0000278'
    'E14' = '  -0.46%  '
    'D15' = '4.091.13'
    'E15' = '  +0.18%  '
    'E16' = '  -0.55%  '
    'D17' = '612.97'
    'E17' = '  -1.08%  '
    'D18' = '3.515.39'
    'E18' = '  +0.03%  '
    'D19' = '70.567.80'
    'E19' = '  +1.01%  '
    'E20' = '  +0.95%  '
    'D21' = '17.56'
    'E21' = '  +1.55%  '
    'E22' = '  -0.08%  '
    'D23' = '9.06'
    'E23' = '  -9.11%  '
    'D24' = '99.72'
    'E24' = '  +3.96%  '
    'D25' = '15.64'
    'E25' = '  -0.60%  '
    'E26' = '  -2.88%  '
    'E27' = '  -0.10%  '
    'E28' = '  -0.78%  '
    'D29' = '34.37'
    'E29' = '  +3.69%  '
    'D30' = '9.06'
    'E30' = '  -1.67%  '
    'E31' = '  -4.10%  '
    'E32' = '  -2.82%  '
    'D33' = '644.66'
    'E33' = '  +13.78%  '
    'B34' = 'NEARProtocol'
    'C34' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D34' = '6.87'
    'E34' = '  -1.40%  '
    'B35' = 'Mantle'
    'C35' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D35' = '1.28'
    'E35' = '  -4.15%  '
    'B36' = 'dogwifhat'
    'C36' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D36' = '3.57'
    'E36' = '  +2.71%  '
    'B37' = 'Hedera'
    'C37' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D37' = '0.0998'
    'E37' = '  -1.32%  '
    'D38' = '10.79'
    'E38' = '  +0.21%  '
    'D39' = '0.0479'
    'E39' = '  +6.63%  '
    'D40' = '56.88'
    'E40' = '  -0.29%  '
    'E41' = '  -0.01%  '
    'E42' = '  +1.45%  '
    'D43' = '0.0₃0749'
    'D44' = '3.371.06'
    'E44' = '  +1.22%  '
    'E45' = '  -4.67%  '
    'D46' = '32.28'
    'E46' = '  -2.42%  '
    'E47' = '  -2.84%  '
    'D48' = '2.57'
    'E48' = '  -2.32%  '
    'E49' = '  +1.00%  '
    'D50' = '133.46'
    'E50' = '  -1.82%  '
    'E51' = '  +0.00%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

